# 2DES BCD AULAS REORGANIZADAS
# Add a new "BCD" attendance column (AZ) mirroring the last PWFE column (AY)
# for every student row, except the hidden row (13).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chamada")

# Row 12 is the only row where the new mark differs ("F" instead of "P").
$faltaRows = @(12)

for ($row = 3; $row -le 30; $row++) {
    if ($row -eq 13) {
        continue
    }

    if ($faltaRows -contains $row) {
        $ws.Cells.Item($row, 52).Value = "F"
    } else {
        $ws.Cells.Item($row, 52).Value = "P"
    }
}

# Restore the active selection to match the reorganised sheet.
$ws.Range("AZ29").Select()
